$wb = $excel.ActiveWorkbook

# Cell updates derived from the canonical OOXML diff.
# Each entry: worksheet name, cell address, new value (or $null to clear the cell).
$updates = @(
    @("ALC", "H13", 8592.333000000001)
    @("ALC", "J13", 8592.333000000001)
    @("ALC", "L13", 8592.333000000001)
    @("ALC", "N13", -8930.333000000001)
    @("ALC", "H20", 20.666666)
    @("ALC", "I20", 20.666666)
    @("ALC", "K20", 20.666666)
    @("ALC", "M20", 209.333334)
    @("ALC", "H35", 20.666666)
    @("ALC", "I35", 20.666666)
    @("ALC", "K35", 20.666666)
    @("ALC", "M35", 358.333334)
    @("ALC", "H40", 1967.375)
    @("ALC", "J40", 2055.5715)
    @("ALC", "L40", 2055.5715)
    @("ALC", "N40", -2405.5715)
    @("ALC", "H45", 38598.5)
    @("ALC", "I45", 7777)
    @("ALC", "K45", 23331)
    @("ALC", "M45", -23139)
    @("ALC", "H52", 2975.8333)
    @("ALC", "I52", 3788.75)
    @("ALC", "J52", 1350)
    @("ALC", "K52", 11366.25)
    @("ALC", "L52", 4050)
    @("ALC", "M52", -11206.25)
    @("ALC", "N52", -4370)
    @("ALC", "H53", 55555904)
    @("ALC", "I53", 148.75)
    @("ALC", "K53", 148.75)
    @("ALC", "M53", 488.25)
    @("ALC", "H98", 3726.0322)
    @("ALC", "J98", 5052.857)
    @("ALC", "L98", 5052.857)
    @("ALC", "N98", -8048.857)
    @("ALC", "H122", 3726.0322)
    @("ALC", "J122", 5052.857)
    @("ALC", "L122", 15158.571)
    @("ALC", "N122", -20058.571)
    @("ALC", "H137", 7693834.5)
    @("ALC", "I137", 1717)
    @("ALC", "J137", 11112553)
    @("ALC", "K137", 5151)
    @("ALC", "L137", 33337659)
    @("ALC", "M137", -2601)
    @("ALC", "N137", -33342759)
    @("ALC", "H141", 1969)
    @("ALC", "I141", 1965.3529)
    @("ALC", "J141", 2000)
    @("ALC", "K141", 5896.0587)
    @("ALC", "L141", 6000)
    @("ALC", "M141", -716.0587000000005)
    @("ALC", "N141", -16360)
    @("ARM", "H32", 289069.8)
    @("ARM", "I32", 358582.97)
    @("ARM", "K32", 358582.97)
    @("ARM", "M32", -358295.97)
    @("ARM", "H122", 2571.6667)
    @("ARM", "I122", 2268.125)
    @("ARM", "K122", 6804.375)
    @("ARM", "M122", -4354.375)
    @("ARM", "H132", 2867.027)
    @("ARM", "I132", 2887.647)
    @("ARM", "J132", 2633.3333)
    @("ARM", "K132", 8662.940999999999)
    @("ARM", "L132", 7899.999899999999)
    @("ARM", "M132", -6132.940999999999)
    @("ARM", "N132", -12959.9999)
    @("BSM", "H134", 24326570)
    @("BSM", "I134", 1887)
    @("BSM", "K134", 5661)
    @("BSM", "M134", -3126)
    @("CRP", "H16", 17861102)
    @("CRP", "I16", 20411260)
    @("CRP", "K16", 20411260)
    @("CRP", "M16", -20410973)
    @("CRP", "H22", 1958.4546)
    @("CRP", "I22", 1943.6666)
    @("CRP", "K22", 1943.6666)
    @("CRP", "M22", -1593.6666)
    @("CRP", "H62", 10712.25)
    @("CRP", "I62", 11159.8)
    @("CRP", "K62", 11159.8)
    @("CRP", "M62", -10535.8)
    @("CRP", "H65", 10712.25)
    @("CRP", "I65", 11159.8)
    @("CRP", "K65", 55799)
    @("CRP", "M65", -52679)
    @("CRP", "H113", 17861102)
    @("CRP", "I113", 20411260)
    @("CRP", "K113", 20411260)
    @("CRP", "M113", -20409090)
    @("CUL", "H62", 14093.111)
    @("CUL", "J62", 15477)
    @("CUL", "L62", 46431)
    @("CUL", "N62", -47803)
    @("CUL", "H63", 10910.244)
    @("CUL", "I63", 0)
    @("CUL", "J63", 10910.244)
    @("CUL", "K63", 0)
    @("CUL", "L63", 32730.732)
    @("CUL", "M63", $null)
    @("CUL", "N63", -34228.732)
    @("CUL", "H64", 9997.6)
    @("CUL", "J64", 9997.6)
    @("CUL", "L64", 29992.8)
    @("CUL", "N64", -30532.8)
    @("CUL", "H65", 14093.111)
    @("CUL", "J65", 15477)
    @("CUL", "L65", 139293)
    @("CUL", "N65", -146157)
    @("CUL", "H66", 10910.244)
    @("CUL", "I66", 0)
    @("CUL", "J66", 10910.244)
    @("CUL", "K66", 0)
    @("CUL", "L66", 98192.19600000001)
    @("CUL", "M66", $null)
    @("CUL", "N66", -105680.196)
    @("CUL", "H67", 9997.6)
    @("CUL", "J67", 9997.6)
    @("CUL", "L67", 29992.8)
    @("CUL", "N67", -31864.8)
    @("CUL", "H92", 1066.6666)
    @("CUL", "I92", 800)
    @("CUL", "J92", 1600)
    @("CUL", "K92", 2400)
    @("CUL", "L92", 4800)
    @("CUL", "M92", -1152)
    @("CUL", "N92", -7296)
    @("GSM", "H2", 95.75)
    @("GSM", "I2", 85.09999999999999)
    @("GSM", "K2", 85.09999999999999)
    @("GSM", "M2", 27.90000000000001)
    @("GSM", "H93", 124777)
    @("GSM", "I93", 0)
    @("GSM", "J93", 124777)
    @("GSM", "K93", 0)
    @("GSM", "L93", 124777)
    @("GSM", "M93", $null)
    @("GSM", "N93", -128521)
    @("GSM", "H95", 21749.75)
    @("GSM", "J95", 21749.75)
    @("GSM", "L95", 21749.75)
    @("GSM", "N95", -27241.75)
    @("GSM", "H102", 100002820)
    @("GSM", "I102", 166668460)
    @("GSM", "K102", 166668460)
    @("GSM", "M102", -166666838)
    @("GSM", "H123", 62499.75)
    @("GSM", "J123", 62499.75)
    @("GSM", "L123", 62499.75)
    @("GSM", "N123", -67399.75)
    @("GSM", "H132", 956034.9399999999)
    @("GSM", "I132", 7815.0557)
    @("GSM", "J132", 2852474.8)
    @("GSM", "K132", 23445.1671)
    @("GSM", "L132", 8557424.399999999)
    @("GSM", "M132", -20915.1671)
    @("GSM", "N132", -8562484.399999999)
    @("LTW", "H22", 7572.6665)
    @("LTW", "I22", 6301)
    @("LTW", "J22", 7663.5)
    @("LTW", "K22", 6301)
    @("LTW", "L22", 7663.5)
    @("LTW", "M22", -6006)
    @("LTW", "N22", -8253.5)
    @("LTW", "H27", 7572.6665)
    @("LTW", "I27", 6301)
    @("LTW", "J27", 7663.5)
    @("LTW", "K27", 6301)
    @("LTW", "L27", 7663.5)
    @("LTW", "M27", -6194)
    @("LTW", "N27", -7877.5)
    @("LTW", "H61", 2278.4688)
    @("LTW", "I61", 2100.5186)
    @("LTW", "K61", 2100.5186)
    @("LTW", "M61", -1898.5186)
    @("LTW", "H68", 1524.7778)
    @("LTW", "I68", 1595)
    @("LTW", "K68", 1595)
    @("LTW", "M68", -846)
    @("LTW", "H71", 1524.7778)
    @("LTW", "I71", 1595)
    @("LTW", "K71", 7975)
    @("LTW", "M71", -4231)
    @("LTW", "H113", 2278.4688)
    @("LTW", "I113", 2100.5186)
    @("LTW", "K113", 2100.5186)
    @("LTW", "M113", 69.48140000000012)
    @("LTW", "H136", 1833.84)
    @("LTW", "I136", 2096.8518)
    @("LTW", "J136", 1736.5616)
    @("LTW", "K136", 6290.555399999999)
    @("LTW", "L136", 5209.6848)
    @("LTW", "M136", -3740.555399999999)
    @("LTW", "N136", -10309.6848)
    @("WVR", "H2", 16332)
    @("WVR", "I2", 19500.5)
    @("WVR", "J2", 9995)
    @("WVR", "K2", 19500.5)
    @("WVR", "L2", 9995)
    @("WVR", "M2", -19388.5)
    @("WVR", "N2", -10219)
    @("WVR", "H4", 22699.4)
    @("WVR", "J4", 20873.75)
    @("WVR", "L4", 20873.75)
    @("WVR", "N4", -21099.75)
    @("WVR", "H132", 2118)
    @("WVR", "I132", 1593.7742)
    @("WVR", "K132", 4781.3226)
    @("WVR", "M132", -2251.3226)
    @("WVR", "H136", 21092.04)
    @("WVR", "I136", 27578.893)
    @("WVR", "K136", 82736.679)
    @("WVR", "M136", -80186.679)
)

foreach ($u in $updates) {
    $sheetName = $u[0]
    $cellAddr = $u[1]
    $newVal = $u[2]
    $ws = $wb.Worksheets.Item($sheetName)
    if ($null -eq $newVal) {
        $ws.Range($cellAddr).ClearContents()
    } else {
        $ws.Range($cellAddr).Value = $newVal
    }
}
